$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find last used row based on column A (Beteckning) starting at row 2
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 125 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value -ne $null) {
        $cell.Value = 45177
    }
}
